$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 230.52
$ws.Range("F30").Value = 29
$ws.Range("G30").Value = 1314.86
$ws.Range("B32").Value = 9009.5
$ws.Range("F51").Value = 63
$ws.Range("G51").Value = 1033.83
$ws.Range("B61").Value = 25376.79
$ws.Range("F97").Value = 58
$ws.Range("G97").Value = 5511.16
$ws.Range("F105").Value = 122
$ws.Range("G105").Value = 5477.8
$ws.Range("F111").Value = 242
$ws.Range("G111").Value = 15415.4
$ws.Range("F116").Value = 151
$ws.Range("G116").Value = 2497.54
$ws.Range("F123").Value = 40
$ws.Range("G123").Value = 8068.4
$ws.Range("B133").Value = 206547.45
$ws.Range("F166").Value = 21
$ws.Range("G166").Value = 1039.08
$ws.Range("F167").Value = 49
$ws.Range("G167").Value = 2073.68
$ws.Range("F174").Value = 22
$ws.Range("G174").Value = 704.4400000000001
$ws.Range("B176").Value = 14898.61
$ws.Range("F292").Value = 20
$ws.Range("G292").Value = 973.8
$ws.Range("B310").Value = 108684.11
$ws.Range("F315").Value = 88
$ws.Range("G315").Value = 9038.48
$ws.Range("F319").Value = 1
$ws.Range("G319").Value = 82.94
$ws.Range("F327").Value = 58
$ws.Range("G327").Value = 7950.06
$ws.Range("F331").Value = 29
$ws.Range("G331").Value = 3435.34
$ws.Range("F351").Value = 3
$ws.Range("G351").Value = 266.4
$ws.Range("F354").Value = 34
$ws.Range("G354").Value = 3437.4
$ws.Range("F372").Value = 93
$ws.Range("G372").Value = 6538.83
$ws.Range("F376").Value = 20
$ws.Range("G376").Value = 1232.8
$ws.Range("B380").Value = 254031.86
$ws.Range("F428").Value = 6
$ws.Range("G428").Value = 400.74
$ws.Range("B429").Value = 6986.72
$ws.Range("B436").Value = 47097
$ws.Range("D436").Value = 112.28
$ws.Range("E436").Value = 134.16
$ws.Range("F436").Value = 15
$ws.Range("G436").Value = 1684.2
$ws.Range("B437").Value = 58047
$ws.Range("D437").Value = 105.54
$ws.Range("E437").Value = 126.1
$ws.Range("F437").Value = 62
$ws.Range("G437").Value = 6543.48
$ws.Range("F456").Value = 5
$ws.Range("G456").Value = 202.7
$ws.Range("B473").Value = 136295.04
$ws.Range("F491").Value = 471
$ws.Range("G491").Value = 6334.95
$ws.Range("F492").Value = 468
$ws.Range("G492").Value = 6154.2
$ws.Range("F493").Value = 547
$ws.Range("G493").Value = 7007.07
$ws.Range("F494").Value = 260
$ws.Range("G494").Value = 6838
$ws.Range("F496").Value = 303
$ws.Range("G496").Value = 4978.29
$ws.Range("F500").Value = 423
$ws.Range("G500").Value = 6861.06
$ws.Range("F502").Value = 927
$ws.Range("G502").Value = 6099.66
$ws.Range("F508").Value = 546
$ws.Range("G508").Value = 8042.58
$ws.Range("B509").Value = 94130.59
$ws.Range("F527").Value = 0
$ws.Range("G527").Value = 0
$ws.Range("B534").Value = 32367.32
$ws.Range("F555").Value = 523
$ws.Range("G555").Value = 3556.4
$ws.Range("B563").Value = 36731.36
$ws.Range("F610").Value = 8
$ws.Range("G610").Value = 27.04
$ws.Range("B613").Value = 6025.98
$ws.Range("F636").Value = 47
$ws.Range("G636").Value = 5773.01
$ws.Range("F639").Value = 15
$ws.Range("G639").Value = 1245.3
$ws.Range("B640").Value = 208294.47
$ws.Range("F663").Value = 2
$ws.Range("G663").Value = 2493.56
$ws.Range("B666").Value = 39984.73
$ws.Range("F668").Value = 9
$ws.Range("G668").Value = 297.99
$ws.Range("F669").Value = 113
$ws.Range("G669").Value = 1792.18
$ws.Range("B677").Value = 20635.01
$ws.Range("F682").Value = 14
$ws.Range("G682").Value = 1148.56
$ws.Range("F683").Value = 24
$ws.Range("G683").Value = 2176.32
$ws.Range("F685").Value = 30
$ws.Range("G685").Value = 1878.3
$ws.Range("F689").Value = 28
$ws.Range("G689").Value = 2660
$ws.Range("F693").Value = 4
$ws.Range("G693").Value = 397.24
$ws.Range("B695").Value = 40698.56
$ws.Range("F702").Value = 45
$ws.Range("G702").Value = 6826.5
$ws.Range("B716").Value = 103573.82
$ws.Range("F761").Value = 29
$ws.Range("G761").Value = 3234.66
$ws.Range("F763").Value = 104
$ws.Range("G763").Value = 2258.88
$ws.Range("F771").Value = 489
$ws.Range("G771").Value = 66019.89
$ws.Range("F772").Value = 17
$ws.Range("G772").Value = 636.14
$ws.Range("F773").Value = 555
$ws.Range("G773").Value = 66994.05
$ws.Range("B775").Value = 246071.38
$ws.Range("F800").Value = 6
$ws.Range("G800").Value = 224.4
$ws.Range("B801").Value = 418.26
$ws.Range("F852").Value = 614
$ws.Range("G852").Value = 18561.22
$ws.Range("F853").Value = 3137
$ws.Range("G853").Value = 511676.07
$ws.Range("F855").Value = 226
$ws.Range("G855").Value = 32690.9
$ws.Range("F859").Value = 180
$ws.Range("G859").Value = 12150
$ws.Range("B861").Value = 626520.23
$ws.Range("B867").Value = 3488051.57
$ws.Range("B868").Value = 3488051.57
